$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 158, shifting existing rows 158:190 down to 159:191
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new record
$ws.Cells.Item(158, 1).Value = 11
$ws.Cells.Item(158, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(158, 3).Value = "Bíobío"
$ws.Cells.Item(158, 4).Value = 45204
$ws.Cells.Item(158, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(158, 5).Value = 8
$ws.Cells.Item(158, 6).Value = 100112001
$ws.Cells.Item(158, 7).Value = "Berenjena"
$ws.Cells.Item(158, 8).Value = "Sin especificar"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 100
$ws.Cells.Item(158, 11).Value = 8500
$ws.Cells.Item(158, 12).Value = 9000
$ws.Cells.Item(158, 13).Value = 8750
$ws.Cells.Item(158, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(158, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(158, 16).Value = 175
$ws.Cells.Item(158, 17).Value = 50
$ws.Cells.Item(158, 18).Value = "Hortaliza"
